# weeklyReportService fix: populate the two empty {fullName} merge-field
# placeholders in the timesheet template (Employee Name field, and the
# "Employee Signature" signature-style line).
#
# Both insertions are done with Range.InsertXML so the resulting markup
# matches what Word itself produces for a merge field typed as
# "{fullName}" -- i.e. three runs ("{", "fullName", "}") with the
# "fullName" word flagged by the spell-checker via <w:proofErr/>.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1) "Employee Name:" row -> second cell, plain {fullName} ----------
$tbl1 = $d.Tables.Item(1)
$nameCell = $tbl1.Cell(1, 2)

$nameFrag = "<w:p $wNs>" +
              "<w:r><w:t>{</w:t></w:r>" +
              '<w:proofErr w:type="spellStart"/>' +
              "<w:r><w:t>fullName</w:t></w:r>" +
              '<w:proofErr w:type="spellEnd"/>' +
              "<w:r><w:t>}</w:t></w:r>" +
            "</w:p>"

$nameCell.Range.InsertXML($nameFrag) | Out-Null

# --- 2) "Employee Signature:" row -> second cell, second (empty) ------
#        paragraph becomes a calligraphy-styled {fullName}, mimicking a
#        signed signature line. The first empty paragraph is untouched.
$tbl3 = $d.Tables.Item(3)
$sigCell = $tbl3.Cell(1, 2)
$sigPara = $sigCell.Range.Paragraphs.Item(2)

$fontProps = '<w:rFonts w:ascii="Lucida Calligraphy" w:eastAsia="Brush Script MT" w:hAnsi="Lucida Calligraphy" w:cs="Apple Chancery"/>'

$sigFrag = "<w:p $wNs>" +
             "<w:pPr><w:rPr>$fontProps</w:rPr></w:pPr>" +
             "<w:r><w:rPr>$fontProps</w:rPr><w:t>{</w:t></w:r>" +
             '<w:proofErr w:type="spellStart"/>' +
             "<w:r><w:rPr>$fontProps<w:i/><w:iCs/></w:rPr><w:t>fullName</w:t></w:r>" +
             '<w:proofErr w:type="spellEnd"/>' +
             "<w:r><w:rPr>$fontProps</w:rPr><w:t>}</w:t></w:r>" +
           "</w:p>"

$sigPara.Range.InsertXML($sigFrag) | Out-Null

Write-Output "done"
